# Initial update to naming conventions (25 GW-high)
#
# Renames the "Baseline-Mid (SC)" sheet to "25 GW - High (SC)".
# Excel automatically propagates the rename to every formula that refers
# to the sheet (including the sheet's own self-referencing CONCAT formulas
# in N2:N26), and updates the active-tab / selection bookkeeping to match
# what was selected interactively when the workbook was saved.

$wb = $excel.ActiveWorkbook

# Rename the sheet - Excel auto-updates all formula references to it.
$oldSheet = $wb.Worksheets.Item("Baseline-Mid (SC)")
$oldSheet.Name = "25 GW - High (SC)"

# The renamed sheet becomes the active tab, with a new selected cell.
$ws = $wb.Worksheets.Item("25 GW - High (SC)")
$ws.Activate()
$ws.Range("I15").Select()
